# fixes indexing bug if window starts at beginning of trace
# Appends the missing trace rows (14-19) that should have been emitted
# when the analysis window starts right at the beginning of a trace.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the existing "Start Time" formatting (column B) as the template for
# the new rows so the new cells share the same time-of-day number format
# already used by the sheet (h:mm:ss), instead of minting a new style.
$timeFormat = $ws.Range("B2").NumberFormat

$newRows = @(
    @(169, 0.45688657407407413, 13),
    @(169, 0.45868055555555554, 9),
    @(178, 0.45682870370370371, 22),
    @(178, 0.45682870370370371, 5),
    @(178, 0.45682870370370371, 14),
    @(178, 0.45682870370370371, 16)
)

$row = 14
foreach ($data in $newRows) {
    $ws.Cells.Item($row, 1).Value = $data[0]

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = $data[1]
    $bCell.NumberFormat = $timeFormat

    $ws.Cells.Item($row, 3).Value = $data[2]

    $row = $row + 1
}

# Move the active selection to where it ends up after entering the data
# (one row below the last new row, column B).
$ws.Range("B20").Select()
